# Weekly update: insert 2 new "Tomate" price records (rows 665-666) ahead of
# the existing history, which shifts down by two rows (665-715 -> 667-717).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 665 - existing rows 665..715 shift to 667..717.
$ws.Rows.Item(665).Insert()
$ws.Rows.Item(665).Insert()

# New row 665
$ws.Cells.Item(665, 1).Value = 11
$ws.Cells.Item(665, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(665, 3).Value = "Bíobío"
$ws.Cells.Item(665, 4).Value = 45021
$ws.Cells.Item(665, 5).Value = 8
$ws.Cells.Item(665, 6).Value = 100112020
$ws.Cells.Item(665, 7).Value = "Tomate"
$ws.Cells.Item(665, 8).Value = "Semiduro"
$ws.Cells.Item(665, 9).Value = "Primera"
$ws.Cells.Item(665, 10).Value = 200
$ws.Cells.Item(665, 11).Value = 9000
$ws.Cells.Item(665, 12).Value = 10000
$ws.Cells.Item(665, 13).Value = 9500
$ws.Cells.Item(665, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(665, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(665, 16).Value = 528
$ws.Cells.Item(665, 17).Value = 18
$ws.Cells.Item(665, 18).Value = "Hortaliza"

# New row 666
$ws.Cells.Item(666, 1).Value = 11
$ws.Cells.Item(666, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(666, 3).Value = "Bíobío"
$ws.Cells.Item(666, 4).Value = 45021
$ws.Cells.Item(666, 5).Value = 8
$ws.Cells.Item(666, 6).Value = 100112020
$ws.Cells.Item(666, 7).Value = "Tomate"
$ws.Cells.Item(666, 8).Value = "Semiduro"
$ws.Cells.Item(666, 9).Value = "Segunda"
$ws.Cells.Item(666, 10).Value = 100
$ws.Cells.Item(666, 11).Value = 8000
$ws.Cells.Item(666, 12).Value = 8000
$ws.Cells.Item(666, 13).Value = 8000
$ws.Cells.Item(666, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(666, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(666, 16).Value = 444
$ws.Cells.Item(666, 17).Value = 18
$ws.Cells.Item(666, 18).Value = "Hortaliza"

$null = $ws.Range("A1").Select()
